$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the paragraph that contains the "Lenovo Yoga 2 Pro" sentence
# together with the CITATION field/content-control that must be
# replaced by a "(bilag XX)" reference.
# ------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Lenovo Yoga 2 Pro*") {
        $targetPara = $p
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not find the paragraph containing 'Lenovo Yoga 2 Pro'"
}

# Capture the paragraph-level formatting (style / numbering / rsids) so
# it can be preserved exactly when the paragraph content is rebuilt.
$openXml = $targetPara.Range.WordOpenXML
$pTag = "<w:p>"
if ($openXml -match '(<w:p [^>]*>)') {
    $pTag = $matches[1]
    # Drop any w14:paraId / w14:textId noise that does not belong to the
    # original document (keep only the classic w:rsid* attributes).
    $pTag = $pTag -replace '\s*w14:paraId="[^"]*"', ''
    $pTag = $pTag -replace '\s*w14:textId="[^"]*"', ''
}
$pPr = ""
if ($openXml -match '(<w:pPr>.*?</w:pPr>)') {
    $pPr = $matches[1]
}

# ------------------------------------------------------------------
# Remove the citation content control (and its field) entirely - this
# is the "Wupti.com, 2014" online-source reference that is being
# replaced with a plain "(bilag XX)" appendix reference.
# ------------------------------------------------------------------
$removed = $false
for ($i = 1; $i -le $d.ContentControls.Count; $i++) {
    $cc = $d.ContentControls.Item($i)
    if ($cc.Range.Text -like "*Wupti.com*") {
        $cc.Delete($true)
        $removed = $true
        break
    }
}
if (-not $removed -and $d.ContentControls.Count -gt 0) {
    $d.ContentControls.Item(1).Delete($true)
}

# ------------------------------------------------------------------
# Rebuild the paragraph's run content:
#   "Lenovo Yoga 2 Pro" + " (bilag XX) " + <bookmark _GoBack> + "benyttes som platform for Fridge app."
# ------------------------------------------------------------------
$newInner = '<w:r><w:t>Lenovo Yoga 2 Pro</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> (bilag XX) </w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
            '<w:r><w:t>benyttes som platform for Fridge app.</w:t></w:r>'

$fragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"' + $pTag.Substring(4)
$fragment = $fragment.Substring(0, $fragment.Length - 1) + '>' + $pPr + $newInner + '</w:p>'

# Re-find the paragraph (its Range may have shifted after the content
# control deletion above) and overwrite its contents in one shot so the
# resulting run/bookmark ordering matches exactly.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Lenovo Yoga 2 Pro*") {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not re-find the paragraph containing 'Lenovo Yoga 2 Pro' after removing the citation"
}

$targetPara.Range.InsertXML($fragment)
